$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the formatting of the last existing data row (row 95) onto the two
# new rows so the new cells pick up the same number formats / styles as
# their neighbours.
$ws.Range("A95:H95").Copy()
$ws.Range("A96:H96").PasteSpecial(-4122)
$ws.Range("A95:H95").Copy()
$ws.Range("A97:H97").PasteSpecial(-4122)

# Row 96: 四方坪站, 2025-10-20 (serial 45948)
$ws.Cells.Item(96, 1).Value = 45948
$ws.Cells.Item(96, 2).Value = "四方坪站"
$ws.Cells.Item(96, 3).Formula = "=18320/126"
$ws.Cells.Item(96, 4).Formula = "=C96/(24*60)"
$ws.Cells.Item(96, 5).Formula = "=10067.02/126"
$ws.Cells.Item(96, 6).Formula = "=3486/126"
$ws.Cells.Item(96, 7).Formula = "=10067.02/(18320/60)"
$ws.Cells.Item(96, 8).Formula = "=434/126"

# Row 97: 高岭站, 2025-10-20 (serial 45948)
$ws.Cells.Item(97, 1).Value = 45948
$ws.Cells.Item(97, 2).Value = "高岭站"
$ws.Cells.Item(97, 3).Formula = "=5801/36"
$ws.Cells.Item(97, 4).Formula = "=C97/(24*60)"
$ws.Cells.Item(97, 5).Formula = "=3944.96/36"
$ws.Cells.Item(97, 6).Formula = "=989.79/36"
$ws.Cells.Item(97, 7).Formula = "=3944.96/(5801/60)"
$ws.Cells.Item(97, 8).Formula = "=149/36"

# Move the active selection the way it ended up after the manual entry.
$ws.Range("J98").Select() | Out-Null
